# Auto-generated Excel COM-interop script
# Applies the 2024-04-09 data update to the crime workbook (350 cell updates across 74 sheets)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 1917
$ws.Range("K3").Value = 1824
$ws.Range("K4").Value = 115
$ws.Range("K5").Value = 1800
$ws.Range("B6").Value = 1697
$ws.Range("K6").Value = 394
$ws.Range("K7").Value = 120
$ws.Range("J8").Value = 29228
$ws.Range("K8").Value = 5878
$ws.Range("K9").Value = 2351
$ws.Range("G10").Value = 41332
$ws.Range("J10").Value = 57339
$ws.Range("K10").Value = 14015
$ws.Range("B11").Value = 104380
$ws.Range("G11").Value = 85343
$ws.Range("K11").Value = 28414

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 45
$ws.Range("K3").Value = 39
$ws.Range("K10").Value = 198
$ws.Range("K11").Value = 481

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("K10").Value = 105
$ws.Range("K11").Value = 152

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K5").Value = 29
$ws.Range("K10").Value = 366
$ws.Range("K11").Value = 596

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 127
$ws.Range("K3").Value = 123
$ws.Range("K6").Value = 23
$ws.Range("K8").Value = 268
$ws.Range("K9").Value = 151
$ws.Range("K10").Value = 375
$ws.Range("K11").Value = 1140

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K2").Value = 11
$ws.Range("K8").Value = 47
$ws.Range("K10").Value = 76
$ws.Range("K11").Value = 178

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 58
$ws.Range("K3").Value = 44
$ws.Range("K5").Value = 19
$ws.Range("K7").Value = 2
$ws.Range("K8").Value = 94
$ws.Range("K10").Value = 132
$ws.Range("K11").Value = 391

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 102
$ws.Range("K8").Value = 142
$ws.Range("K9").Value = 70
$ws.Range("K10").Value = 232
$ws.Range("K11").Value = 691

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K3").Value = 61
$ws.Range("K5").Value = 37
$ws.Range("K9").Value = 38
$ws.Range("K10").Value = 227
$ws.Range("K11").Value = 570

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 273
$ws.Range("K4").Value = 139
$ws.Range("K6").Value = 275
$ws.Range("K7").Value = 620
$ws.Range("K8").Value = 1140
$ws.Range("K9").Value = 108
$ws.Range("K10").Value = 316
$ws.Range("K11").Value = 567
$ws.Range("K12").Value = 115
$ws.Range("K15").Value = 262
$ws.Range("K16").Value = 178
$ws.Range("K18").Value = 196
$ws.Range("K19").Value = 670
$ws.Range("K20").Value = 487
$ws.Range("K22").Value = 108
$ws.Range("K23").Value = 284
$ws.Range("K24").Value = 157
$ws.Range("K26").Value = 58
$ws.Range("K27").Value = 416
$ws.Range("K29").Value = 875
$ws.Range("K31").Value = 277
$ws.Range("K32").Value = 53
$ws.Range("K33").Value = 691
$ws.Range("K35").Value = 73
$ws.Range("K36").Value = 373
$ws.Range("K37").Value = 702
$ws.Range("K39").Value = 43
$ws.Range("K41").Value = 133
$ws.Range("K42").Value = 723
$ws.Range("K43").Value = 314
$ws.Range("K44").Value = 357
$ws.Range("K47").Value = 245
$ws.Range("J48").Value = 3095
$ws.Range("K48").Value = 738
$ws.Range("K49").Value = 466
$ws.Range("K50").Value = 257
$ws.Range("K51").Value = 352
$ws.Range("K52").Value = 481
$ws.Range("K53").Value = 596
$ws.Range("K54").Value = 988
$ws.Range("K57").Value = 126
$ws.Range("K59").Value = 72
$ws.Range("K60").Value = 178
$ws.Range("K61").Value = 34
$ws.Range("B63").Value = 1495
$ws.Range("G63").Value = 1864
$ws.Range("J63").Value = 492
$ws.Range("K63").Value = 146
$ws.Range("K64").Value = 251
$ws.Range("K65").Value = 389
$ws.Range("K66").Value = 154
$ws.Range("J67").Value = 2748
$ws.Range("K67").Value = 558
$ws.Range("K69").Value = 152
$ws.Range("K72").Value = 169
$ws.Range("K73").Value = 370
$ws.Range("K74").Value = 83
$ws.Range("K75").Value = 116
$ws.Range("K76").Value = 697
$ws.Range("J77").Value = 524
$ws.Range("K77").Value = 108
$ws.Range("K78").Value = 583
$ws.Range("K79").Value = 570
$ws.Range("K82").Value = 78
$ws.Range("K83").Value = 391
$ws.Range("K84").Value = 204
$ws.Range("K85").Value = 985
$ws.Range("K86").Value = 223
$ws.Range("J88").Value = 1063
$ws.Range("K88").Value = 253
$ws.Range("K89").Value = 492
$ws.Range("K90").Value = 290
$ws.Range("K91").Value = 238
$ws.Range("K92").Value = 88
$ws.Range("K93").Value = 205
$ws.Range("K94").Value = 701
$ws.Range("K96").Value = 475
$ws.Range("K97").Value = 471
$ws.Range("K98").Value = 357
$ws.Range("K100").Value = 79
$ws.Range("B101").Value = 104380
$ws.Range("G101").Value = 85343
$ws.Range("K101").Value = 28414

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("K2").Value = 10
$ws.Range("K8").Value = 24
$ws.Range("K11").Value = 116

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K9").Value = 14
$ws.Range("K11").Value = 251

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 55
$ws.Range("K3").Value = 72
$ws.Range("K4").Value = 10
$ws.Range("K8").Value = 166
$ws.Range("K9").Value = 74
$ws.Range("K11").Value = 702

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K3").Value = 14
$ws.Range("K8").Value = 82
$ws.Range("K10").Value = 238
$ws.Range("K11").Value = 416

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 45
$ws.Range("K9").Value = 68
$ws.Range("K10").Value = 131
$ws.Range("K11").Value = 389

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K2").Value = 27
$ws.Range("K8").Value = 72
$ws.Range("K11").Value = 277

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 73
$ws.Range("K3").Value = 80
$ws.Range("K8").Value = 109
$ws.Range("K9").Value = 84
$ws.Range("J10").Value = 886
$ws.Range("K10").Value = 158
$ws.Range("J11").Value = 2748
$ws.Range("K11").Value = 558

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K10").Value = 97
$ws.Range("K11").Value = 204

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K5").Value = 25
$ws.Range("K8").Value = 97
$ws.Range("K10").Value = 503
$ws.Range("K11").Value = 701

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K2").Value = 18
$ws.Range("K8").Value = 82
$ws.Range("K9").Value = 55
$ws.Range("K10").Value = 486
$ws.Range("K11").Value = 697

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("K8").Value = 8
$ws.Range("K10").Value = 55
$ws.Range("K11").Value = 79

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("K10").Value = 122
$ws.Range("K11").Value = 178

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K10").Value = 342
$ws.Range("K11").Value = 466

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K9").Value = 39
$ws.Range("K10").Value = 251
$ws.Range("K11").Value = 471

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K2").Value = 31
$ws.Range("K10").Value = 197
$ws.Range("K11").Value = 475

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 38
$ws.Range("K8").Value = 137
$ws.Range("K9").Value = 44
$ws.Range("K10").Value = 714
$ws.Range("K11").Value = 988

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K10").Value = 198
$ws.Range("K11").Value = 370

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 85
$ws.Range("K3").Value = 109
$ws.Range("K8").Value = 219
$ws.Range("K9").Value = 108
$ws.Range("K10").Value = 281
$ws.Range("K11").Value = 875

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K3").Value = 13
$ws.Range("K5").Value = 68
$ws.Range("K6").Value = 13
$ws.Range("K8").Value = 91
$ws.Range("K9").Value = 36
$ws.Range("J10").Value = 2124
$ws.Range("K10").Value = 501
$ws.Range("J11").Value = 3095
$ws.Range("K11").Value = 738

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K3").Value = 54
$ws.Range("K8").Value = 172
$ws.Range("K9").Value = 61
$ws.Range("K10").Value = 257
$ws.Range("K11").Value = 670

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("K8").Value = 21
$ws.Range("K10").Value = 87
$ws.Range("K11").Value = 154

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K2").Value = 11
$ws.Range("K8").Value = 70
$ws.Range("K9").Value = 29
$ws.Range("K10").Value = 179
$ws.Range("K11").Value = 357

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K3").Value = 62
$ws.Range("K5").Value = 40
$ws.Range("K8").Value = 173
$ws.Range("K9").Value = 101
$ws.Range("K10").Value = 275
$ws.Range("K11").Value = 723

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("K2").Value = 8
$ws.Range("K8").Value = 44
$ws.Range("K11").Value = 108

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K3").Value = 56
$ws.Range("K8").Value = 161
$ws.Range("K10").Value = 228
$ws.Range("K11").Value = 620

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K2").Value = 19
$ws.Range("K5").Value = 24
$ws.Range("K10").Value = 114
$ws.Range("K11").Value = 275

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K10").Value = 40
$ws.Range("K11").Value = 133

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K5").Value = 31
$ws.Range("K8").Value = 116
$ws.Range("K10").Value = 146
$ws.Range("K11").Value = 373

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K10").Value = 207
$ws.Range("K11").Value = 316

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K10").Value = 148
$ws.Range("K11").Value = 223

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K3").Value = 20
$ws.Range("K8").Value = 128
$ws.Range("K10").Value = 334
$ws.Range("K11").Value = 583

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("K5").Value = 20
$ws.Range("K10").Value = 82
$ws.Range("K11").Value = 157

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value = 19
$ws.Range("K10").Value = 126
$ws.Range("K11").Value = 262

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 44
$ws.Range("K8").Value = 116
$ws.Range("K10").Value = 284
$ws.Range("K11").Value = 567

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K5").Value = 9
$ws.Range("K9").Value = 15
$ws.Range("K10").Value = 140
$ws.Range("K11").Value = 284

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K10").Value = 89
$ws.Range("K11").Value = 238

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K8").Value = 70
$ws.Range("K10").Value = 184
$ws.Range("K11").Value = 352

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 43
$ws.Range("K5").Value = 24
$ws.Range("K8").Value = 124
$ws.Range("K11").Value = 487

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K10").Value = 124
$ws.Range("K11").Value = 245

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 26
$ws.Range("K10").Value = 131
$ws.Range("K11").Value = 290

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("K10").Value = 144
$ws.Range("K11").Value = 257

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("K2").Value = 9
$ws.Range("K5").Value = 14
$ws.Range("K8").Value = 73
$ws.Range("K11").Value = 205

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K5").Value = 13
$ws.Range("K10").Value = 78
$ws.Range("K11").Value = 196

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K2").Value = 23
$ws.Range("J10").Value = 152
$ws.Range("J11").Value = 524
$ws.Range("K11").Value = 108

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K9").Value = 21
$ws.Range("K10").Value = 148
$ws.Range("K11").Value = 273

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K3").Value = 6
$ws.Range("K10").Value = 116
$ws.Range("K11").Value = 169

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K3").Value = 19
$ws.Range("K10").Value = 163
$ws.Range("K11").Value = 314

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("K10").Value = 61
$ws.Range("K11").Value = 139

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K10").Value = 308
$ws.Range("K11").Value = 492

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("K10").Value = 59
$ws.Range("K11").Value = 73

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("K8").Value = 17
$ws.Range("K10").Value = 22
$ws.Range("K11").Value = 53

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K8").Value = 40
$ws.Range("K9").Value = 34
$ws.Range("K10").Value = 245
$ws.Range("K11").Value = 357

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("K8").Value = 8
$ws.Range("K11").Value = 58

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("K10").Value = 33
$ws.Range("K11").Value = 88

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range("K10").Value = 14
$ws.Range("K11").Value = 34

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("K9").Value = 50
$ws.Range("K10").Value = 78

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K8").Value = 30
$ws.Range("K11").Value = 108

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 120
$ws.Range("K5").Value = 111
$ws.Range("K8").Value = 175
$ws.Range("K10").Value = 363
$ws.Range("K11").Value = 985

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("K10").Value = 35
$ws.Range("K11").Value = 72

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("K8").Value = 24
$ws.Range("K9").Value = 43

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K10").Value = 66
$ws.Range("K11").Value = 126

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K3").Value = 17
$ws.Range("K9").Value = 49
$ws.Range("J10").Value = 429
$ws.Range("K10").Value = 106
$ws.Range("J11").Value = 1063
$ws.Range("K11").Value = 253

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("K10").Value = 70
$ws.Range("K11").Value = 83

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("K10").Value = 79
$ws.Range("K11").Value = 115
